$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking text (e.g. '1.00', '0.170') that
# Excel would otherwise auto-coerce to a Number on assignment, stripping
# formatting / trailing zeros. Force Text format first so the literal
# string is preserved exactly, matching the source data feed's formatting.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.129.44"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.883.98"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "526.33"
$ws.Range("E5").Value = "  +8.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.48"
$ws.Range("E6").Value = "  -2.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.608"
$ws.Range("E7").Value = "  -1.97%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.717"
$ws.Range("E9").Value = "  -2.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.170"
$ws.Range("E10").Value = "  +1.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000331"
$ws.Range("E11").Value = "  -4.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "41.83"
$ws.Range("E12").Value = "  -2.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.516.81"
$ws.Range("E13").Value = "  -0.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.16"
$ws.Range("E14").Value = "  -5.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.042.92"
$ws.Range("E15").Value = "  +3.85%  "
$ws.Range("E16").Value = "  +7.85%  "
$ws.Range("E17").Value = "  -0.66%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.72"
$ws.Range("E18").Value = "  -3.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.63"
$ws.Range("E19").Value = "  -2.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.111.15"
$ws.Range("E20").Value = "  +1.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "423.33"
$ws.Range("E21").Value = "  -1.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.33"
$ws.Range("E22").Value = "  -4.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.11"
$ws.Range("E23").Value = "  -7.00%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.58"
$ws.Range("E24").Value = "  -1.05%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.05"
$ws.Range("E25").Value = "  +9.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.65"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.51"
$ws.Range("E27").Value = "  -5.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "35.92"
$ws.Range("E28").Value = "  -4.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "694.70"
$ws.Range("E29").Value = "  -2.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "13.09"
$ws.Range("E30").Value = "  -5.11%  "
$ws.Range("E31").Value = "  -4.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.81"
$ws.Range("E32").Value = "  -3.73%  "
$ws.Range("E33").Value = "  +11.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.441"
$ws.Range("E34").Value = "  +10.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.92"
$ws.Range("E35").Value = "  -4.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "40.15"
$ws.Range("E36").Value = "  -3.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0833"
$ws.Range("E37").Value = "  -7.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.149"
$ws.Range("E38").Value = "  +3.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.998"
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0477"
$ws.Range("E41").Value = "  -2.22%  "
$ws.Range("E42").Value = "  -7.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.01"
$ws.Range("E43").Value = "  +0.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.95"
$ws.Range("E44").Value = "  -5.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.32"
$ws.Range("E45").Value = "  -0.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.139"
$ws.Range("E46").Value = "  -1.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.04"
$ws.Range("E47").Value = "  +8.34%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₆0341"
$ws.Range("E48").Value = "  +2.49%  "
$ws.Range("B49").Value = "LidoDAOToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.27"
$ws.Range("E49").Value = "  -4.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "142.39"
$ws.Range("E50").Value = "  -1.84%  "
$ws.Range("E51").Value = "  -4.07%  "
